$wb = $excel.ActiveWorkbook

# =========================================================================
# Generate Report for Handoff
# Adds a new handed-off file (cabc622a-053b-4896-9b25-fe977a147a9d.md) as
# row 3 to the Overview sheet and to both locale detail sheets (zh-cn,
# de-de), mirroring the existing row 2 pattern for 18414fa7-....md.
# =========================================================================

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = "cabc622a-053b-4896-9b25-fe977a147a9d.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 08:50:57"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ebea501f9cf2634bf39845e3adbe5f84ca7b39eb/e2e/cabc622a-053b-4896-9b25-fe977a147a9d.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\cabc622a-053b-4896-9b25-fe977a147a9d.md") | Out-Null

# --- zh-cn sheet --------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = "cabc622a-053b-4896-9b25-fe977a147a9d.00bb567247de2f11bc028f95b8e2ae60525e2e24.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-09-01 08:50:52"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("O3").Value = "'False"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ebea501f9cf2634bf39845e3adbe5f84ca7b39eb/e2e/cabc622a-053b-4896-9b25-fe977a147a9d.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "cabc622a-053b-4896-9b25-fe977a147a9d.md") | Out-Null

# --- de-de sheet ----------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = "cabc622a-053b-4896-9b25-fe977a147a9d.00bb567247de2f11bc028f95b8e2ae60525e2e24.de-de.xlf"
$wsDe.Range("H3").Value = "2016-09-01 08:50:57"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("O3").Value = "'False"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ebea501f9cf2634bf39845e3adbe5f84ca7b39eb/e2e/cabc622a-053b-4896-9b25-fe977a147a9d.md", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "cabc622a-053b-4896-9b25-fe977a147a9d.md") | Out-Null

Write-Host "Handoff report row added to Overview, zh-cn, de-de"
